# Add six new color-palette rows to the "經典流派與風格" sheet (4th tab),
# update its selection/scroll position + page setup, and update the
# selection on the "現代生活與情感" sheet (3rd tab).

$wb = $excel.ActiveWorkbook

# --- Sheet 4: 經典流派與風格 -------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Clone the formatting (font/border/alignment) of the last existing data
# row (row 7) down into the six new rows (8-13) before filling in values.
$ws4.Range("A7:D7").Copy()
$ws4.Range("A8:D13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 - 森林晨曦 (Olive & Moss)
$ws4.Range("A8").Value = "森林晨曦 (Olive & Moss)"
$ws4.Range("B8").Value = "#506736"
$ws4.Range("C8").Value = "#91A263"
$ws4.Range("D8").Value = "自然橄欖與苔蘚的柔和過渡"
$ws4.Rows.Item(8).RowHeight = 83.25

# Row 9 - 寒地松林 (Pine & Bone)
$ws4.Range("A9").Value = "寒地松林 (Pine & Bone)"
$ws4.Range("B9").Value = "#2C341B"
$ws4.Range("C9").Value = "#DDDFC2"
$ws4.Range("D9").Value = "沉穩松針與骨色的冷調質感"
$ws4.Rows.Item(9).RowHeight = 83.25

# Row 10 - 復古大地 (Mustard & Camel)
$ws4.Range("A10").Value = "復古大地 (Mustard & Camel)"
$ws4.Range("B10").Value = "#756633"
$ws4.Range("C10").Value = "#BC9E5F"
$ws4.Range("D10").Value = "芥末綠與駝色的復古大地感"
$ws4.Rows.Item(10).RowHeight = 83.25

# Row 11 - 迷彩灰綠 (Feldgrau & Bud)
$ws4.Range("A11").Value = "迷彩灰綠 (Feldgrau & Bud)"
$ws4.Range("B11").Value = "#396153"
$ws4.Range("C11").Value = "#E4E8B8"
$ws4.Range("D11").Value = "軍事風格與淡色春芽的撞色"
$ws4.Rows.Item(11).RowHeight = 83.25

# Row 12 - 深空月石 (Sparkle & Moonstone)
$ws4.Range("A12").Value = "深空月石 (Sparkle & Moonstone)"
$ws4.Range("B12").Value = "#406768"
$ws4.Range("C12").Value = "#6FA9BB"
$ws4.Range("D12").Value = "深邃藍綠與月石色的冷冽星空"
$ws4.Rows.Item(12).RowHeight = 99.75

# Row 13 - 經典商務 (Prussian & Rhus)
$ws4.Range("A13").Value = "經典商務 (Prussian & Rhus)"
$ws4.Range("B13").Value = "#003153"
$ws4.Range("C13").Value = "#E3A841"
$ws4.Range("D13").Value = "櫨色與深藍的經典沉穩配比"
$ws4.Rows.Item(13).RowHeight = 83.25

# Page setup for sheet 4
$ws4.PageSetup.PaperSize = 9
$ws4.PageSetup.Orientation = 1

# Scroll / selection for sheet 4: show the new rows, with B8:D13 selected
$ws4.Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range("B8:D13").Select()

# --- Sheet 3: 現代生活與情感 --------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1:D3").Select()

# Leave sheet 4 as the active tab, matching the workbook's saved state.
$ws4.Select()
